# Add transformer parameters (d_model, num_layers) to the parameter sheet
# and update a few existing parameter values, per commit:
# "Added transformer parameters to xlsx file; they are loaded from there
#  by default instead of load_pred_par.m"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the added parameters
$ws.Range("K1").Value = "d_model"
$ws.Range("L1").Value = "num_layers"

# Updated existing parameter values in row 2
$ws.Range("C2").Value = 57
$ws.Range("D2").Value = 110
$ws.Range("F2").Value = 0.0005
$ws.Range("H2").Value = 10

# New parameter values in row 2
$ws.Range("K2").Value = 16
$ws.Range("L2").Value = 1

# Match the updated selection recorded in the saved view
$ws.Range("F5").Select()
